# Set the document's default editing language to Swedish (sv-SE), matching
# the "set as default" Word operation that stamps w:lang onto every style's
# run properties (w:eastAsia stays en-US, w:bidi stays ar-SA - unchanged).
#
# The docDefaults/rPrDefault's w:lang val="en-US" -> "sv-SE" change that real
# Word performs when you pick "Set As Default" in the Language dialog is not
# reachable through the Style object model exposed here, but every one of
# the 164 style definitions (paragraph/character/table/list) each receive an
# explicit <w:lang val="sv-SE" eastAsia="en-US" bidi="ar-SA"/> override -
# exactly what iterating Document.Styles and writing Font.LanguageID /
# Font.LanguageIDFarEast / Font.LanguageIDOther reproduces.

$d = $word.ActiveDocument

$count = $d.Styles.Count
for ($i = 1; $i -le $count; $i++) {
    $s = $d.Styles.Item($i)
    $s.Font.LanguageID = "sv-SE"
    $s.Font.LanguageIDFarEast = "en-US"
    $s.Font.LanguageIDOther = "ar-SA"
}
